$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "If there is an issue in generating recommendations*") {
        $p.Range.Delete()
        break
    }
}
